$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 16:35"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 1531465
$ws.Range("C4").Value = 3801
$ws.Range("E4").Value = 1094011
$ws.Range("G4").Value = 82
$ws.Range("H4").Value = 91060

# Row 11: 'Alemania' -> 'Alemania'
$ws.Range("B11").Value = 176933
$ws.Range("C11").Value = 282
$ws.Range("E11").Value = 14271
$ws.Range("G11").Value = 13
$ws.Range("H11").Value = 8062

# Row 14: 'India' -> 'India'
$ws.Range("B14").Value = 97577
$ws.Range("C14").Value = 1879
$ws.Range("E14").Value = 57712

# Row 62: 'Moldavia' -> 'Moldavia'
$ws.Range("B62").Value = 6138
$ws.Range("C62").Value = 78
$ws.Range("D62").Value = 2425
$ws.Range("E62").Value = 3496
$ws.Range("G62").Value = 6
$ws.Range("H62").Value = 217

# Row 75: 'Uzbekistan' -> 'Uzbekistan'
$ws.Range("D75").Value = 2314
$ws.Range("E75").Value = 452

# Row 108: 'Republica de Chipre' -> 'Libano'
$ws.Range("A108").Value = "Libano"
$ws.Range("B108").Value = 931
$ws.Range("C108").Value = 20
$ws.Range("D108").Value = 247
$ws.Range("E108").Value = 658
$ws.Range("H108").Value = 26

# Row 109: 'Kenia' -> 'Republica de Chipre'
$ws.Range("A109").Value = "Republica de Chipre"
$ws.Range("B109").Value = 916
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 515
$ws.Range("E109").Value = 384
$ws.Range("H109").Value = 17

# Row 110: 'Libano' -> 'Kenia'
$ws.Range("A110").Value = "Kenia"
$ws.Range("B110").Value = 912
$ws.Range("C110").Value = 25
$ws.Range("D110").Value = 336
$ws.Range("E110").Value = 526
$ws.Range("H110").Value = 50

# Row 179: 'Macao' -> 'Macao'
$ws.Range("D179").Value = 44
$ws.Range("E179").Value = 1

# Row 195: 'Santa Lucia' -> 'Nueva Caledonia'
$ws.Range("A195").Value = "Nueva Caledonia"

# Row 197: 'Nueva Caledonia' -> 'Santa Lucia'
$ws.Range("A197").Value = "Santa Lucia"

# Row 208: 'Groenlandia' -> 'Surinam'
$ws.Range("A208").Value = "Surinam"
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 9
$ws.Range("E208").Value = 1
$ws.Range("H208").Value = 1

# Row 209: 'Montserrat' -> 'Seychelles'
$ws.Range("A209").Value = "Seychelles"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

# Row 210: 'Seychelles' -> 'Groenlandia'
$ws.Range("A210").Value = "Groenlandia"

# Row 211: 'Surinam' -> 'Montserrat'
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("D211").Value = 10

# Row 214: 'San Bartolome' -> 'Sahara Occidental'
$ws.Range("A214").Value = "Sahara Occidental"

# Row 215: 'Sahara Occidental' -> 'San Bartolome'
$ws.Range("A215").Value = "San Bartolome"
